# Applies the "Updated cryptos list" price/volume refresh described in the commit.
# Column D ("Price") values that look numeric (e.g. "72.46") must stay plain text,
# exactly like the rest of the sheet (all cells are inline strings). Assigning a bare
# numeric-looking string lets Excel auto-convert the cell to a real number, so each
# such value is written with a leading apostrophe (a literal `'` character placed
# in front of the text -- in a PowerShell single-quoted string that's written as two
# consecutive quote characters, `''`, followed by the value and the closing quote).
# This mirrors typing '72.46 into the cell in the Excel UI and forces a text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''64.445.58'
$ws.Range("E2").Value = '  -0.19%  '

# Row 3
$ws.Range("D3").Value = '''3.413.60'
$ws.Range("E3").Value = '  -1.25%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '''572.01'
$ws.Range("E5").Value = '  -0.51%  '

# Row 6
$ws.Range("D6").Value = '''156.70'
$ws.Range("E6").Value = '  -3.00%  '

# Row 7
$ws.Range("D7").Value = '''0.616'
$ws.Range("E7").Value = '  +6.81%  '

# Row 9
$ws.Range("D9").Value = '''3.416.93'
$ws.Range("E9").Value = '  -1.19%  '

# Row 10
$ws.Range("D10").Value = '''7.18'
$ws.Range("E10").Value = '  -2.26%  '

# Row 11
$ws.Range("D11").Value = '''0.122'
$ws.Range("E11").Value = '  -2.56%  '

# Row 12
$ws.Range("D12").Value = '''0.441'
$ws.Range("E12").Value = '  +0.27%  '

# Row 13
$ws.Range("D13").Value = '''4.005.33'

# Row 14
$ws.Range("E14").Value = '  +0.10%  '

# Row 15
$ws.Range("E15").Value = '  -3.17%  '

# Row 16
$ws.Range("D16").Value = '''27.99'
$ws.Range("E16").Value = '  -2.77%  '

# Row 17
$ws.Range("D17").Value = '''64.482.99'
$ws.Range("E17").Value = '  -0.13%  '

# Row 18
$ws.Range("D18").Value = '''3.423.05'
$ws.Range("E18").Value = '  -0.80%  '

# Row 19
$ws.Range("E19").Value = '  -0.39%  '

# Row 20
$ws.Range("D20").Value = '''13.97'
$ws.Range("E20").Value = '  -2.54%  '

# Row 21
$ws.Range("D21").Value = '''375.44'
$ws.Range("E21").Value = '  -3.75%  '

# Row 22
$ws.Range("D22").Value = '''7.99'
$ws.Range("E22").Value = '  -2.17%  '

# Row 23
$ws.Range("D23").Value = '''0.554'
$ws.Range("E23").Value = '  +1.19%  '

# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''72.46'
$ws.Range("E24").Value = '  -0.56%  '

# Row 25
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''0.998'
$ws.Range("E25").Value = '  +0.04%  '

# Row 26
$ws.Range("E26").Value = '  -4.01%  '

# Row 27
$ws.Range("D27").Value = '''10.25'
$ws.Range("E27").Value = '  +8.16%  '

# Row 28
$ws.Range("E28").Value = '  -2.38%  '

# Row 29
$ws.Range("E29").Value = '  +0.08%  '

# Row 30
$ws.Range("D30").Value = '''1.50'
$ws.Range("E30").Value = '  +4.31%  '

# Row 31
$ws.Range("E31").Value = '  -0.38%  '

# Row 33
$ws.Range("D33").Value = '''23.08'
$ws.Range("E33").Value = '  -2.21%  '

# Row 34
$ws.Range("E34").Value = '  +1.77%  '

# Row 35
$ws.Range("D35").Value = '''1.61'
$ws.Range("E35").Value = '  +6.83%  '

# Row 36
$ws.Range("D36").Value = '''160.47'
$ws.Range("E36").Value = '  -0.63%  '

# Row 37
$ws.Range("D37").Value = '''1.89'
$ws.Range("E37").Value = '  -0.52%  '

# Row 38
$ws.Range("E38").Value = '  -1.48%  '

# Row 39
$ws.Range("E39").Value = '  +5.19%  '

# Row 40
$ws.Range("D40").Value = '''26.62'
$ws.Range("E40").Value = '  -2.98%  '

# Row 41
$ws.Range("D41").Value = '''2.844.35'
$ws.Range("E41").Value = '  -2.81%  '

# Row 42
$ws.Range("D42").Value = '''4.61'
$ws.Range("E42").Value = '  +1.24%  '

# Row 43
$ws.Range("E43").Value = '  +0.25%  '

# Row 44
$ws.Range("D44").Value = '''26.49'
$ws.Range("E44").Value = '  +9.99%  '

# Row 45
$ws.Range("E45").Value = '  -0.94%  '

# Row 46
$ws.Range("D46").Value = '''0.768'
$ws.Range("E46").Value = '  -0.79%  '

# Row 47
$ws.Range("D47").Value = '''319.54'
$ws.Range("E47").Value = '  +7.64%  '

# Row 48
$ws.Range("E48").Value = '  -0.89%  '

# Row 49
$ws.Range("E49").Value = '  +1.93%  '

# Row 50
$ws.Range("D50").Value = '''6.59'
$ws.Range("E50").Value = '  +0.93%  '

# Row 51
$ws.Range("D51").Value = '''0.854'
$ws.Range("E51").Value = '  -2.21%  '
